# Auto-generated edit script: updates crypto price/volume table to reflect
# the latest GitHub Actions scrape (commit "Updated cryptos list ...").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure these Price cells stay plain text (they look numeric but the
# source data models them as text, e.g. "214.10" must not become 214.1)
$textCells = @('D5', 'D6', 'D8', 'D9', 'D10', 'D11', 'D12', 'D14', 'D16', 'D19', 'D20', 'D23', 'D24', 'D25', 'D26', 'D27', 'D29', 'D31', 'D33', 'D35', 'D38', 'D41', 'D42', 'D44', 'D46', 'D47', 'D48', 'D49', 'D50', 'D51')
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '25.911.93'
$ws.Range("E2").Value = '  -0.15%  '

$ws.Range("D3").Value = '1.636.05'
$ws.Range("E3").Value = '  -0.39%  '

$ws.Range("E4").Value = '  -0.28%  '

$ws.Range("D5").Value = '214.10'
$ws.Range("E5").Value = '  -0.69%  '

$ws.Range("D6").Value = '0.5061'
$ws.Range("E6").Value = '  -0.45%  '

$ws.Range("E7").Value = '  -0.30%  '

$ws.Range("D8").Value = '0.2569'
$ws.Range("E8").Value = '  +0.26%  '

$ws.Range("D9").Value = '0.06355'
$ws.Range("E9").Value = '  -0.48%  '

$ws.Range("D10").Value = '19.67'
$ws.Range("E10").Value = '  +0.96%  '

$ws.Range("D11").Value = '0.07738'
$ws.Range("E11").Value = '  -0.72%  '

$ws.Range("D12").Value = '4.278'

$ws.Range("D13").Value = '1.640.86'
$ws.Range("E13").Value = '  -0.04%  '

$ws.Range("D14").Value = '0.5436'
$ws.Range("E14").Value = '  -0.66%  '

$ws.Range("D15").Value = '0.0₅7731'
$ws.Range("E15").Value = '  -1.56%  '

$ws.Range("D16").Value = '64.04'
$ws.Range("E16").Value = '  -0.63%  '

$ws.Range("D17").Value = '25.927.88'
$ws.Range("E17").Value = '  -0.35%  '

$ws.Range("E18").Value = '  -0.21%  '

$ws.Range("B19").Value = 'BitcoinCash'
$ws.Range("C19").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D19").Value = '195.40'
$ws.Range("E19").Value = '  -1.42%  '

$ws.Range("B20").Value = 'Uniswap'
$ws.Range("C20").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D20").Value = '4.430'
$ws.Range("E20").Value = '  -0.46%  '

$ws.Range("E21").Value = '  -0.71%  '

$ws.Range("E22").Value = '  +0.73%  '

$ws.Range("D23").Value = '1.002'
$ws.Range("E23").Value = '  -0.57%  '

$ws.Range("D24").Value = '1.891'
$ws.Range("E24").Value = '  +0.81%  '

$ws.Range("D25").Value = '142.98'
$ws.Range("E25").Value = '  +1.40%  '

$ws.Range("D26").Value = '0.1244'
$ws.Range("E26").Value = '  +8.28%  '

$ws.Range("D27").Value = '6.816'
$ws.Range("E27").Value = '  -1.00%  '

$ws.Range("E28").Value = '  -0.90%  '

$ws.Range("D29").Value = '1.236'
$ws.Range("E29").Value = '  -0.43%  '

$ws.Range("E30").Value = '  -3.17%  '

$ws.Range("D31").Value = '3.239'
$ws.Range("E31").Value = '  -0.71%  '

$ws.Range("E32").Value = '  +0.00%  '

$ws.Range("D33").Value = '1.548'
$ws.Range("E33").Value = '  +0.23%  '

$ws.Range("E34").Value = '  +0.11%  '

$ws.Range("D35").Value = '0.9105'
$ws.Range("E35").Value = '  +1.19%  '

$ws.Range("E36").Value = '  -0.90%  '

$ws.Range("B37").Value = 'Maker'
$ws.Range("C37").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D37").Value = '1.125.40'
$ws.Range("E37").Value = '  -0.69%  '

$ws.Range("B38").Value = 'ImmutableX'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D38").Value = '0.5495'
$ws.Range("E38").Value = '  -0.16%  '

$ws.Range("E40").Value = '  -0.48%  '

$ws.Range("D41").Value = '5.577'
$ws.Range("E41").Value = '  -0.79%  '

$ws.Range("D42").Value = '0.8039'
$ws.Range("E42").Value = '  -1.82%  '

$ws.Range("E43").Value = '  -8.52%  '

$ws.Range("D44").Value = '98.58'
$ws.Range("E44").Value = '  -1.68%  '

$ws.Range("D45").Value = '1.768.73'
$ws.Range("E45").Value = '  -0.63%  '

$ws.Range("D46").Value = '0.4480'
$ws.Range("E46").Value = '  -1.10%  '

$ws.Range("D47").Value = '1.003'
$ws.Range("E47").Value = '  -0.02%  '

$ws.Range("D48").Value = '54.98'
$ws.Range("E48").Value = '  -0.01%  '

$ws.Range("D49").Value = '0.05159'
$ws.Range("E49").Value = '  +1.72%  '

$ws.Range("D50").Value = '7.522'
$ws.Range("E50").Value = '  +1.98%  '

$ws.Range("D51").Value = '0.9987'
$ws.Range("E51").Value = '  -0.77%  '

